# This edit re-shuffles the per-record weekly price data (columns D, J, K,
# L, M, P for rows 2..39) across the existing rows, i.e. each target row
# ends up with the Fecha/Volumen/Precio* values that used to belong to a
# different (source) row. All other columns (A, B, C, E-I, N, O, Q, R)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: targetRow -> sourceRow (the row whose current D/J/K/L/M/P values
# should be written into targetRow).
$rowMap = @{
    2 = 25
    3 = 10
    4 = 28
    5 = 6
    6 = 37
    7 = 26
    8 = 20
    9 = 24
    10 = 34
    11 = 32
    12 = 7
    13 = 39
    14 = 15
    15 = 12
    16 = 18
    17 = 17
    18 = 2
    19 = 36
    20 = 11
    21 = 33
    22 = 30
    23 = 35
    24 = 22
    25 = 8
    26 = 14
    27 = 29
    28 = 5
    29 = 19
    30 = 9
    31 = 13
    32 = 4
    33 = 38
    34 = 3
    35 = 27
    36 = 21
    37 = 31
    38 = 23
    39 = 16
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot the current (pre-edit) values for every affected cell before we
# start overwriting anything, since several rows read from rows that will
# themselves be overwritten later in the loop.
$snapshot = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 39; $r++) {
        $addr = "$col$r"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    if ($sourceRow -eq $targetRow) {
        continue
    }
    foreach ($col in $cols) {
        $srcAddr = "$col$sourceRow"
        $dstAddr = "$col$targetRow"
        $ws.Range($dstAddr).Value2 = $snapshot[$srcAddr]
    }
}
